# Update fermentation and separation improvement parameters in uncertainty scenarios
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fermentation TAL yield baseline (row 23): 0.68 -> 0.73
$ws.Range("E23").Value = 0.73

# Fermentation TAL titer baseline (row 24): 76 -> 68
$ws.Range("E24").Value = 68

# Reflect the selection change recorded in the saved view state
$ws.Range("A23:XFD24").Select()

$wb.Save()
